$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BF column holds a "Date" label in row 1 and date strings in rows 2-31.
# The date strings were recorded one day off (due to how NBA stats were
# displayed) and need to be corrected from "5-2-2012-13" to "2013-05-02".
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    if ($cell.Value -eq "5-2-2012-13") {
        $cell.Value = "2013-05-02"
    }
}
